# Add a "Discussed" marker down column X (Curation status helper column)
# for every data row of the "Source all" sheet, and update the active
# selection to match (mirrors the author's online-editing-support commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source all")
$ws.Activate()

# Column X already carries the "Curation status" header in X1; fill the
# data rows (2 through 141) with the new "Discussed" value.
$ws.Range("X2:X141").Value = "Discussed"

# Leave the selection on the newly-filled range, as in the edited workbook.
$ws.Range("X2:X141").Select()
